$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list (prices / 1h volume change %) per latest scrape.
# For Price cells (column D) whose new value looks like a plain number
# (e.g. "215.40", "21.34"), force the cell to Text format first so Excel
# keeps the exact literal string (incl. trailing zeros) instead of
# re-interpreting it as a numeric value.

$ws.Range('D2').Value = '27.075.71'
$ws.Range('E2').Value = '  +0.62%  '
$ws.Range('D3').Value = '1.676.09'
$ws.Range('E3').Value = '  +0.43%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.40'
$ws.Range('E5').Value = '  +0.41%  '
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.256'
$ws.Range('E8').Value = '  +2.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '21.34'
$ws.Range('E9').Value = '  +5.40%  '
$ws.Range('E10').Value = '  +0.22%  '
$ws.Range('D12').Value = '1.912.67'
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('D13').Value = '1.670.51'
$ws.Range('E13').Value = '  +0.28%  '
$ws.Range('E14').Value = '  +1.15%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.538'
$ws.Range('E15').Value = '  +2.23%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.16'
$ws.Range('E16').Value = '  +1.04%  '
$ws.Range('D17').Value = '27.052.20'
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '237.55'
$ws.Range('E18').Value = '  +1.84%  '
$ws.Range('E19').Value = '  +2.06%  '
$ws.Range('D20').Value = '0.0₃0739'
$ws.Range('E20').Value = '  +0.78%  '
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('E22').Value = '  +1.48%  '
$ws.Range('E23').Value = '  +1.95%  '
$ws.Range('E24').Value = '  -1.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.44'
$ws.Range('E25').Value = '  +0.76%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.24'
$ws.Range('E26').Value = '  +1.98%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.41'
$ws.Range('E27').Value = '  +3.10%  '
$ws.Range('E28').Value = '  +0.90%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  +0.33%  '
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('E32').Value = '  +0.83%  '
$ws.Range('D33').Value = '1.532.01'
$ws.Range('E33').Value = '  +5.34%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.19'
$ws.Range('E34').Value = '  +2.18%  '
$ws.Range('E35').Value = '  +2.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.39'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.594'
$ws.Range('E37').Value = '  +1.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.919'
$ws.Range('E38').Value = '  +2.56%  '
$ws.Range('E39').Value = '  +2.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.08'
$ws.Range('E40').Value = '  +2.53%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '67.77'
$ws.Range('E42').Value = '  +2.41%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.51'
$ws.Range('E43').Value = '  -3.81%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.26'
$ws.Range('E44').Value = '  -1.40%  '
$ws.Range('D45').Value = '1.820.20'
$ws.Range('E46').Value = '  +0.57%  '
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('E48').Value = '  +1.15%  '
$ws.Range('E49').Value = '  +2.49%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.04'
$ws.Range('E50').Value = '  +5.91%  '
$ws.Range('E51').Value = '  +0.51%  '
